# The source change (commit "Moving from POI 3.17.0 to 4.0.1") only
# re-serialized the expected-generation fixture with a newer Apache POI
# version: every hunk in the diff is a pure XML attribute/namespace
# reordering (e.g. <w:pgSz w:h="16838" w:w="11906"/> -> <w:pgSz w:w="11906"
# w:h="16838"/>, xmlns:* declaration order on the root elements, etc.).
# No text, formatting, style value, or structural content actually differs
# between the "before" and "after" XML - every attribute/element value is
# identical, only the on-disk attribute order changed.
#
# The Word object model does not expose a way to control raw XML attribute
# or namespace-declaration ordering (that's purely a side effect of which
# library/version wrote the package), so there is no document-level edit to
# perform here. Touch the active document without changing anything.
$d = $word.ActiveDocument
$d.Content | Out-Null
